$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Update the vote counters at the top of the sheet
$ws.Range("B2").Value = 8
$ws.Range("B3").Value = 8
$ws.Range("B4").Value = 16

# Update the selected cell to match the new selection
$ws.Range("J8").Select() | Out-Null

# Tableau1 (bleu votes) - columns B (Prenom) / C (nom), rows 8-15
$ws.Range("B8").Value  = "Anatol"
$ws.Range("C8").Value  = "Anatol"

$ws.Range("B9").Value  = "Annabel"
$ws.Range("C9").Value  = "Arange"

$ws.Range("B10").Value = "billiab"
$ws.Range("C10").Value = "Anatol"

$ws.Range("B11").Value = "Anatol"
$ws.Range("C11").Value = "Anatol"

$ws.Range("B12").Value = "Bill"
$ws.Range("C12").Value = "Bill"

$ws.Range("B13").Value = "chali"
$ws.Range("C13").Value = "Bill"

$ws.Range("B14").Value = "Anatol"
$ws.Range("C14").Value = "Barnon"

$ws.Range("B15").Value = "Anatol"
$ws.Range("C15").Value = "Anatol"

# Tableau2 (rouge votes) - columns E (Prenom) / F (nom), rows 8-15
$ws.Range("E8").Value  = "Anatol"
$ws.Range("F8").Value  = "Anatol"

$ws.Range("E9").Value  = "Annabel"
$ws.Range("F9").Value  = "Arange"

$ws.Range("E10").Value = "billiab"
$ws.Range("F10").Value = "Anatol"

$ws.Range("E11").Value = "Anatol"
$ws.Range("F11").Value = "Anatol"

$ws.Range("E12").Value = "Bill"
$ws.Range("F12").Value = "Bill"

$ws.Range("E13").Value = "chali"
$ws.Range("F13").Value = "Bill"

$ws.Range("E14").Value = "Anatol"
$ws.Range("F14").Value = "Barnon"

$ws.Range("E15").Value = "Anatol"
$ws.Range("F15").Value = "Anatol"
